# 25. Find all elements that appear more than " n/k " times.
#
# Fill in the missing "Sl. No." numbering in column B for the Binary Search
# Tree / Greedy sections (rows 214-272). The numbering continues from the
# last filled value right above this block (B211 = 196), skipping the blank
# separator rows that have no data in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$counter = 196
for ($row = 214; $row -le 272; $row++) {
    $topic = $ws.Cells.Item($row, 1).Text
    if ($topic -ne "") {
        $counter = $counter + 1
        $ws.Cells.Item($row, 2).Value = $counter
    }
}

# Scroll the sheet down to the newly-edited area and leave the selection on
# the cell right after the last edited row, mirroring the saved view state.
$ws.Range("C273").Select()
$excel.ActiveWindow.ScrollRow = 266
$excel.ActiveWindow.ScrollColumn = 1
